$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = "xgbclassifier"
$ws.Range("B5").Value = 0.8295648464163823
$ws.Range("C5").Value = 0.8149205055034652
$ws.Range("D5").Value = 0.8528156996587031
$ws.Range("E5").Value = 0.8334375651448822
$ws.Range("F5").Value = 0.6598435045992631
$ws.Range("G5").Value = 0.6591296928327646
$ws.Range("H5").Value = 0.8295648464163823
